$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F3").Value = 1751
$ws.Range("F4").Value = 167
$ws.Range("F5").Value = 471
$ws.Range("F6").Value = 830
$ws.Range("F7").Value = 256
$ws.Range("F8").Value = 1220
$ws.Range("F9").Value = 350
$ws.Range("F11").Value = 884
$ws.Range("F12").Value = 703
$ws.Range("F13").Value = 191
$ws.Range("F15").Value = 144
$ws.Range("F18").Value = 2954
$ws.Range("F19").Value = 2629
$ws.Range("F24").Value = 233
$ws.Range("F26").Value = 5317
$ws.Range("F28").Value = 988
$ws.Range("F29").Value = 25
$ws.Range("F31").Value = 323
$ws.Range("F32").Value = 1109
$ws.Range("F33").Value = 70
$ws.Range("F34").Value = 57
$ws.Range("F35").Value = 291
$ws.Range("F36").Value = 37

$ws = $wb.Worksheets.Item(2)
$ws.Range("F4").Value = 1134
$ws.Range("F7").Value = 232
$ws.Range("F17").Value = 988
$ws.Range("F23").Value = 4
$ws.Range("F26").Value = 3956
$ws.Range("F34").Value = 33

$ws = $wb.Worksheets.Item(3)
$ws.Range("F3").Value = 78
$ws.Range("F5").Value = 2472
$ws.Range("F6").Value = 1052
$ws.Range("F9").Value = 1336
$ws.Range("F10").Value = 367

$ws = $wb.Worksheets.Item(4)
$ws.Range("F4").Value = 2472
$ws.Range("F5").Value = 1751
$ws.Range("F6").Value = 1052
$ws.Range("F7").Value = 1336
$ws.Range("F8").Value = 367
$ws.Range("F10").Value = 167
$ws.Range("F11").Value = 471
$ws.Range("F12").Value = 830
$ws.Range("F13").Value = 256
$ws.Range("F14").Value = 1220
$ws.Range("F15").Value = 350
$ws.Range("F16").Value = 884
$ws.Range("F17").Value = 703
$ws.Range("F18").Value = 1134
$ws.Range("F19").Value = 1134
$ws.Range("F20").Value = 191
$ws.Range("F23").Value = 2954
$ws.Range("F24").Value = 2629
$ws.Range("F28").Value = 233
$ws.Range("F29").Value = 5317
$ws.Range("F31").Value = 988
$ws.Range("F34").Value = 25
$ws.Range("F37").Value = 323
$ws.Range("F46").Value = 1109
$ws.Range("F50").Value = 57
$ws.Range("F51").Value = 291
$ws.Range("F52").Value = 37
